# Update the "Lương" (salary) sheet with the new computed figures after
# adding the combined salary report logic.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B11").Value = -0
$ws.Range("B12").Value = 9.5
$ws.Range("B13").Value = 1357142.857142857
$ws.Range("B20").Value = -1000000
$ws.Range("B30").Value = -0
$ws.Range("B32").Value = 457142.8571428573
$ws.Range("B34").Value = 457142.8571428573
